$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B3 from "Gernert" to "Gruschka"
$ws.Range("B3").Value = "Gruschka"

# Delete the now-obsolete last row (previously row 4: 3, Gruschka)
$ws.Rows("4").Delete()

# Adjust the window size of the workbook view
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Width = 5040
$win.Height = 10245
